$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("C65").Value = "'" + '=(List)$Step2["Calc"]'
$ws1.Range("C66").Value = "'" + '=((Map)$Step3[0])["Step1"]'
$ws1.Range("C67").Value = "'" + '=((Map)$Step3[1])["Step1"]'
$ws1.Range("C93").Value = "'" + '=$Step2["Calc"]'
$ws1.Range("C95").Value = "'" + '=((Passport)((MyType)((Map)$Step4)["Step7"]).passportData).passportId'
$ws1.Range("C96").Value = "'" + '=((MyType)((Map)$Step4)["Step7"]).someMap'
$ws1.Range("C97").Value = "'" + '=(List)((MyType)((Map)$Step4)["Step7"]).someList'
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("C13").Value = "'" + '=(Map[][])$Step2["Calc"]'
$ws2.Range("C18").Value = "'" + '=((MyType)((Map)$Step3[0][0])["Step7"]).someText'
$ws2.Range("C19").Value = "'" + '=((MyType)((Map)$Step3[0][1])["Step7"]).someText'
$ws2.Range("C20").Value = "'" + '=((MyType)((Map)$Step3[1][0])["Step7"]).someText'
$ws2.Range("C21").Value = "'" + '=((MyType)((Map)$Step3[1][1])["Step7"]).someText'
$ws2.Range("F26").Value = "'" + '_res_.$Step4["Step1"]:Integer'
$ws2.Range("G26").Value = "'" + '_res_.$Step5["Step1"]:Integer'
$ws2.Range("H26").Value = "'" + '_res_.$Step6["Step1"]:Integer'
$ws2.Range("I26").Value = "'" + '_res_.$Step7["Step1"]:Integer'
$ws2.Range("F27").Value = "'" + '_res_.$Step4["Step1"]:Integer'
$ws2.Range("G27").Value = "'" + '_res_.$Step5["Step1"]:Integer'
$ws2.Range("H27").Value = "'" + '_res_.$Step6["Step1"]:Integer'
$ws2.Range("I27").Value = "'" + '_res_.$Step7["Step1"]:Integer'
$ws2.Range("C45").Value = "'" + '=(List)$Step2["Calc"]'
$ws2.Range("C50").Value = "'" + '=((MyType)(((Map[])$Step3[0])[0])["Step7"]).someText'
$ws2.Range("C51").Value = "'" + '=((MyType)(((Map[])$Step3[0])[1])["Step7"]).someText'
$ws2.Range("C52").Value = "'" + '=((MyType)(((Map[])$Step3[1])[0])["Step7"]).someText'
$ws2.Range("C53").Value = "'" + '=((MyType)(((Map[])$Step3[1])[1])["Step7"]).someText'
$ws2.Range("F58").Value = "'" + '_res_.$Step4["Step1"]:Integer'
$ws2.Range("G58").Value = "'" + '_res_.$Step5["Step1"]:Integer'
$ws2.Range("H58").Value = "'" + '_res_.$Step6["Step1"]:Integer'
$ws2.Range("I58").Value = "'" + '_res_.$Step7["Step1"]:Integer'
$ws2.Range("F59").Value = "'" + '_res_.$Step4["Step1"]:Integer'
$ws2.Range("G59").Value = "'" + '_res_.$Step5["Step1"]:Integer'
$ws2.Range("H59").Value = "'" + '_res_.$Step6["Step1"]:Integer'
$ws2.Range("I59").Value = "'" + '_res_.$Step7["Step1"]:Integer'
$ws2.Range("C75").Value = "'" + '=(List)$Step1.toMap()["Calc"]'
$ws2.Range("C76").Value = "'" + '=(Map)((List)(((Map)$Step2[0])["Calc"]))[0]'
$ws2.Range("C77").Value = "'" + '=(Map)((List)(((Map)$Step2[0])["Calc"]))[1]'
$ws2.Range("C78").Value = "'" + '=(Map)((List)(((Map)$Step2[1])["Calc"]))[0]'
$ws2.Range("C79").Value = "'" + '=(Map)((List)(((Map)$Step2[1])["Calc"]))[1]'
$ws2.Range("C80").Value = "'" + '=((MyType)(((Map)((List)(((Map)$Step2[1])["Calc"]))[1])["Step7"])).someMap'
$ws2.Range("D83").Value = "'" + '_res_.$Step3["Step1"]:Integer'
$ws2.Range("E83").Value = "'" + '_res_.$Step4["Step1"]:Integer'
$ws2.Range("F83").Value = "'" + '_res_.$Step5["Step1"]:Integer'
$ws2.Range("G83").Value = "'" + '_res_.$Step6["Step1"]:Integer'
$ws2.Range("D84").Value = "'" + '_res_.$Step3["Step1"]'
$ws2.Range("E84").Value = "'" + '_res_.$Step4["Step1"]'
$ws2.Range("F84").Value = "'" + '_res_.$Step5["Step1"]'
$ws2.Range("G84").Value = "'" + '_res_.$Step6["Step1"]'
